# Fix a recurring typo "Tiếp nhân" -> "Tiếp nhận" in the "Ràng buộc" (constraint)
# column (column C) of Sheet1. This affects rows 6, 7 and 8, which share two
# distinct strings that both contain the misspelling.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C6").Value = "Chỉ user thuộc bộ phận Tiếp nhận&Phóng thích mới được phép sử dụng"
$ws.Range("C7").Value = "Chỉ user thuộc Bộ phận Tiếp nhận&Phóng thích và Bộ phận Quản lý mới được phép sử dụng"
$ws.Range("C8").Value = "Chỉ user thuộc bộ phận Tiếp nhận&Phóng thích mới được phép sử dụng"

# Reflect the active cell/selection that was saved with the workbook.
$ws.Activate()
$ws.Range("B4").Select() | Out-Null
